# "oracle init database first version done"
#
# The Database sheet previously reused a single shared admin account
# (gameAdmin2 / admin00!!) for all three environment rows. This splits
# each row out to its own admin user + password:
#   row2 (test_game_system) -> gameAdmin1 / admin01
#   row3 (test_game_logDB)  -> gameAdmin2 / admin02   (user unchanged)
#   row4 (test_game_zoneDB) -> gameAdmin3 / admin03

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Database" sheet is already the active/selected tab

# -- User / Password columns (C = User, D = Password) --
$ws.Range("C2").Value = "gameAdmin1"
$ws.Range("C4").Value = "gameAdmin3"

$ws.Range("D2").Value = "admin01"
$ws.Range("D3").Value = "admin02"
$ws.Range("D4").Value = "admin03"

# -- widen the User column (no longer auto-fit) to fit the longer values --
$ws.Columns.Item(3).ColumnWidth = 13.36

# -- cursor/selection left on D7 when the file was saved --
$ws.Range("D7").Select() | Out-Null
